$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing the existing row 2 (Colombian Primera B)
# down to row 3 and the existing row 3 (Brazilian Serie A) down to row 4.
$ws.Rows.Item(2).Insert()
# The inserted row inherits formatting from the row above it (the bold header);
# strip that back to the plain/default style used by the other data rows.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the Bosnian Premier League match.
$ws.Range("A2").Value = "Bosnian Premier League"

# "2025-11-12" looks like a date to Excel's auto-detection, so force it to be
# stored as plain text (matching the other Date cells in the sheet) and then
# drop the temporary text number-format back to the default style.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-11-12"
$ws.Range("B2").ClearFormats()

$ws.Range("C2").Value = "14:00:00"
$ws.Range("D2").Value = "Borac Banja Luka"
$ws.Range("E2").Value = "Zrinjski"
$ws.Range("F2").Value = 1.04
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.04
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.03
$ws.Range("K2").Value = 950
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.26
$ws.Range("O2").Value = 1.01
$ws.Range("P2").Value = 1.26
$ws.Range("Q2").Value = 1.02
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 1.43
$ws.Range("T2").Value = 1.04
$ws.Range("U2").Value = 1.04
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# The Colombian Primera B row (now shifted down to row 3) also had some of
# its odds updated: Odd_H_Back and Odd_A_Back from 1.04 to 1.33, and
# Odd_D_Lay from 1000 to 4.
$ws.Range("F3").Value = 1.33
$ws.Range("H3").Value = 1.33
$ws.Range("K3").Value = 4
